# Update "想去人数" (interest count) figures that were refreshed by the
# automated gh-pages data generator.
#
# Sheet "展览" (exhibitions) and sheet "全部类型" (all types) both list the
# same events in rows 2, 7, 11, 12 and 16 - column F holds the interest
# count for each event.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 2;  Value = 730 },
    @{ Row = 7;  Value = 78 },
    @{ Row = 11; Value = 4798 },
    @{ Row = 12; Value = 4544 },
    @{ Row = 16; Value = 33 }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($update in $updates) {
        $ws.Cells.Item($update.Row, 6).Value = $update.Value
    }
}
